$d = $word.ActiveDocument

# Replace the text of the first run ("O vídeo for") with the new commit text.
$d.Content.Find.Execute("O vídeo for", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Mudança para subir no git", 2)

# Remove the remainder of the first paragraph's old text (the second run after the bookmark).
$d.Content.Find.Execute("nece uma maneira poderosa de ajudá-lo a provar seu argumento. Ao clicar em Vídeo Online, você pode colar o código de inserção do vídeo que deseja adicionar. Você também pode digitar uma palavra-chave para pesquisar online o vídeo mais adequado ao seu documento.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "", 2)

# Delete the second paragraph entirely ("Para dar ao documento...").
$p2 = $d.Paragraphs(2)
$p2.Range.Text = ""
$p2.Range.Delete()

# After deleting paragraph 2, the former paragraph 3 ("Temas e estilos...") is now paragraph 2.
$p3 = $d.Paragraphs(2)
$p3.Range.Text = ""
$p3.Range.Delete()

# After deleting that, the trailing empty paragraph is now paragraph 2; remove it too.
$p4 = $d.Paragraphs(2)
$p4.Range.Text = ""
$p4.Range.Delete()
